$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 0.06515966666666667
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 10.95398031243089
$ws.Range("R2").Value = 98.585822811878
$ws.Range("S2").Value = 0.2984182258032519
$ws.Range("T2").Value = 0.298418225803252

# Row 3
$ws.Range("G3").Value = 0.06515966666666667
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("R3").Value = 95.59288782565301
$ws.Range("S3").Value = 0.2893586437755394
$ws.Range("T3").Value = 0.2893586437755394

# Row 4
$ws.Range("G4").Value = 0.06515966666666667
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 10.81608308362333
$ws.Range("R4").Value = 97.34474775261
$ws.Range("S4").Value = 0.294661504941043
$ws.Range("T4").Value = 0.294661504941043

# Row 5
$ws.Range("G5").Value = 0.06515966666666667
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 4.315311933581889
$ws.Range("R5").Value = 38.837807402237
$ws.Range("S5").Value = 0.1175616254801657
$ws.Range("T5").Value = 0.1175616254801657
